$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.033695101737976
$ws.Range("B1").Value = 3.088480949401855
$ws.Range("C1").Value = 6.727529525756836
$ws.Range("D1").Value = 1.864564061164856
$ws.Range("E1").Value = 1.316844940185547
